$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of row 2 (A2:O2) while keeping formatting/style intact
$ws.Range("A2:O2").ClearContents()

# Update the active selection cell as recorded in the saved workbook
$ws.Range("T28").Select()
